# Phantom_Profits workbook update: rebases currentAveragePrice/LevePrice/LeveProfit
# figures to a newer market snapshot, and strips the bold/centered/thin-bordered
# header style (A1:N1) back to the workbook default on every sheet.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Sheet: ALC
# ---------------------------------------------------------------------------
$ws = $wb.Worksheets.Item("ALC")

# Reset header row (A1:N1) to the default style: removes bold font,
# thin box border, and center/top alignment.
$ws.Range("A1:N1").ClearFormats()

$ws.Cells.Item(20, 8).Value = 632.2857  # H20: 679.3333 -> 632.2857
$ws.Cells.Item(20, 9).Value = 671  # I20: 735.2 -> 671
$ws.Cells.Item(20, 11).Value = 671  # K20: 735.2 -> 671
$ws.Cells.Item(20, 13).Value = -441  # M20: -505.2 -> -441

$ws.Cells.Item(28, 8).Value = 1655.3  # H28: 1705.2 -> 1655.3
$ws.Cells.Item(28, 9).Value = 1693.5  # I28: 1727.3334 -> 1693.5
$ws.Cells.Item(28, 10).Value = 1502.5  # J28: 1506 -> 1502.5
$ws.Cells.Item(28, 11).Value = 1693.5  # K28: 1727.3334 -> 1693.5
$ws.Cells.Item(28, 12).Value = 1502.5  # L28: 1506 -> 1502.5
$ws.Cells.Item(28, 13).Value = -1208.5  # M28: -1242.3334 -> -1208.5
$ws.Cells.Item(28, 14).Value = -2472.5  # N28: -2476 -> -2472.5

$ws.Cells.Item(34, 8).Value = 6962  # H34: 7720.25 -> 6962
$ws.Cells.Item(34, 9).Value = 6962  # I34: 7720.25 -> 6962
$ws.Cells.Item(34, 11).Value = 6962  # K34: 7720.25 -> 6962
$ws.Cells.Item(34, 13).Value = -6759  # M34: -7517.25 -> -6759

$ws.Cells.Item(35, 8).Value = 632.2857  # H35: 679.3333 -> 632.2857
$ws.Cells.Item(35, 9).Value = 671  # I35: 735.2 -> 671
$ws.Cells.Item(35, 11).Value = 671  # K35: 735.2 -> 671
$ws.Cells.Item(35, 13).Value = -292  # M35: -356.2 -> -292

$ws.Cells.Item(36, 8).Value = 6962  # H36: 7720.25 -> 6962
$ws.Cells.Item(36, 9).Value = 6962  # I36: 7720.25 -> 6962
$ws.Cells.Item(36, 11).Value = 6962  # K36: 7720.25 -> 6962
$ws.Cells.Item(36, 13).Value = -6247  # M36: -7005.25 -> -6247

$ws.Cells.Item(43, 8).Value = 5158.8184  # H43: 5812.375 -> 5158.8184
$ws.Cells.Item(43, 9).Value = 5218.375  # I43: 5714.143 -> 5218.375
$ws.Cells.Item(43, 10).Value = 5000  # J43: 6500 -> 5000
$ws.Cells.Item(43, 11).Value = 5218.375  # K43: 5714.143 -> 5218.375
$ws.Cells.Item(43, 12).Value = 5000  # L43: 6500 -> 5000
$ws.Cells.Item(43, 13).Value = -5149.375  # M43: -5645.143 -> -5149.375
$ws.Cells.Item(43, 14).Value = -5138  # N43: -6638 -> -5138

$ws.Cells.Item(53, 8).Value = 474.82352  # H53: 458.1111 -> 474.82352
$ws.Cells.Item(53, 9).Value = 433.81818  # I53: 412.16666 -> 433.81818
$ws.Cells.Item(53, 11).Value = 433.81818  # K53: 412.16666 -> 433.81818
$ws.Cells.Item(53, 13).Value = 203.18182  # M53: 224.83334 -> 203.18182

$ws.Cells.Item(82, 8).Value = 1414.2  # H82: 3313.6667 -> 1414.2
$ws.Cells.Item(82, 9).Value = 1414.2  # I82: 2997 -> 1414.2
$ws.Cells.Item(82, 10).Value = 0  # J82: 3947 -> 0
$ws.Cells.Item(82, 11).Value = 4242.6  # K82: 8991 -> 4242.6
$ws.Cells.Item(82, 12).Value = 0  # L82: 11841 -> 0
$ws.Cells.Item(82, 13).Value = -3836.6  # M82: -8585 -> -3836.6
$ws.Cells.Item(82, 14).Value = $null  # N82: -12653 -> (removed)

$ws.Cells.Item(85, 8).Value = 1414.2  # H85: 3313.6667 -> 1414.2
$ws.Cells.Item(85, 9).Value = 1414.2  # I85: 2997 -> 1414.2
$ws.Cells.Item(85, 10).Value = 0  # J85: 3947 -> 0
$ws.Cells.Item(85, 11).Value = 4242.6  # K85: 8991 -> 4242.6
$ws.Cells.Item(85, 12).Value = 0  # L85: 11841 -> 0
$ws.Cells.Item(85, 13).Value = -2838.6  # M85: -7587 -> -2838.6
$ws.Cells.Item(85, 14).Value = $null  # N85: -14649 -> (removed)

$ws.Cells.Item(111, 8).Value = 2814  # H111: 3027.111 -> 2814
$ws.Cells.Item(111, 9).Value = 1948.75  # I111: 2299.6667 -> 1948.75
$ws.Cells.Item(111, 11).Value = 5846.25  # K111: 6899.000100000001 -> 5846.25
$ws.Cells.Item(111, 13).Value = -2779.25  # M111: -3832.000100000001 -> -2779.25

$ws.Cells.Item(125, 8).Value = 2200  # H125: 1766.6666 -> 2200
$ws.Cells.Item(125, 9).Value = 0  # I125: 900 -> 0
$ws.Cells.Item(125, 11).Value = 0  # K125: 8100 -> 0
$ws.Cells.Item(125, 13).Value = $null  # M125: -5640 -> (removed)

$ws.Cells.Item(135, 8).Value = 662.3  # H135: 687.3333 -> 662.3
$ws.Cells.Item(135, 9).Value = 655.2857  # I135: 691.6667 -> 655.2857
$ws.Cells.Item(135, 11).Value = 5897.571300000001  # K135: 6225.0003 -> 5897.571300000001
$ws.Cells.Item(135, 13).Value = -3362.571300000001  # M135: -3690.0003 -> -3362.571300000001

# ---------------------------------------------------------------------------
# Sheet: ARM
# ---------------------------------------------------------------------------
$ws = $wb.Worksheets.Item("ARM")

# Reset header row (A1:N1) to the default style: removes bold font,
# thin box border, and center/top alignment.
$ws.Range("A1:N1").ClearFormats()

$ws.Cells.Item(45, 8).Value = 1886.5714  # H45: 1775.625 -> 1886.5714
$ws.Cells.Item(45, 9).Value = 1886.5714  # I45: 1775.625 -> 1886.5714
$ws.Cells.Item(45, 11).Value = 1886.5714  # K45: 1775.625 -> 1886.5714
$ws.Cells.Item(45, 13).Value = -1509.5714  # M45: -1398.625 -> -1509.5714

$ws.Cells.Item(110, 8).Value = 9051.182  # H110: 9787.2 -> 9051.182
$ws.Cells.Item(110, 9).Value = 9456.4  # I110: 10319.223 -> 9456.4
$ws.Cells.Item(110, 11).Value = 9456.4  # K110: 10319.223 -> 9456.4
$ws.Cells.Item(110, 13).Value = -7411.4  # M110: -8274.223 -> -7411.4

$ws.Cells.Item(122, 8).Value = 458.16666  # H122: 529.4 -> 458.16666
$ws.Cells.Item(122, 9).Value = 458.16666  # I122: 529.4 -> 458.16666
$ws.Cells.Item(122, 11).Value = 1374.49998  # K122: 1588.2 -> 1374.49998
$ws.Cells.Item(122, 13).Value = 1075.50002  # M122: 861.8000000000002 -> 1075.50002

$ws.Cells.Item(130, 8).Value = 37330  # H130: 39496.25 -> 37330
$ws.Cells.Item(130, 10).Value = 37330  # J130: 39496.25 -> 37330
$ws.Cells.Item(130, 12).Value = 37330  # L130: 39496.25 -> 37330
$ws.Cells.Item(130, 14).Value = -47370  # N130: -49536.25 -> -47370

$ws.Cells.Item(131, 8).Value = 79975  # H131: 79976 -> 79975
$ws.Cells.Item(131, 10).Value = 79975  # J131: 79976 -> 79975
$ws.Cells.Item(131, 12).Value = 79975  # L131: 79976 -> 79975
$ws.Cells.Item(131, 14).Value = -90055  # N131: -90056 -> -90055

# ---------------------------------------------------------------------------
# Sheet: BSM
# ---------------------------------------------------------------------------
$ws = $wb.Worksheets.Item("BSM")

# Reset header row (A1:N1) to the default style: removes bold font,
# thin box border, and center/top alignment.
$ws.Range("A1:N1").ClearFormats()

$ws.Cells.Item(43, 8).Value = 475000  # H43: 0 -> 475000
$ws.Cells.Item(43, 10).Value = 475000  # J43: 0 -> 475000
$ws.Cells.Item(43, 12).Value = 475000  # L43: 0 -> 475000
$ws.Cells.Item(43, 14).Value = -475362  # N43: None -> -475362

$ws.Cells.Item(86, 8).Value = 803.7143  # H86: 1020 -> 803.7143
$ws.Cells.Item(86, 9).Value = 604.3333  # I86: 693.3333 -> 604.3333
$ws.Cells.Item(86, 11).Value = 604.3333  # K86: 693.3333 -> 604.3333
$ws.Cells.Item(86, 13).Value = 518.6667  # M86: 429.6667 -> 518.6667

$ws.Cells.Item(89, 8).Value = 803.7143  # H89: 1020 -> 803.7143
$ws.Cells.Item(89, 9).Value = 604.3333  # I89: 693.3333 -> 604.3333
$ws.Cells.Item(89, 11).Value = 3021.6665  # K89: 3466.6665 -> 3021.6665
$ws.Cells.Item(89, 13).Value = 2594.3335  # M89: 2149.3335 -> 2594.3335

$ws.Cells.Item(115, 8).Value = 0  # H115: 80500 -> 0
$ws.Cells.Item(115, 10).Value = 0  # J115: 80500 -> 0
$ws.Cells.Item(115, 12).Value = 0  # L115: 80500 -> 0
$ws.Cells.Item(115, 14).Value = $null  # N115: -83634 -> (removed)

# ---------------------------------------------------------------------------
# Sheet: CRP
# ---------------------------------------------------------------------------
$ws = $wb.Worksheets.Item("CRP")

# Reset header row (A1:N1) to the default style: removes bold font,
# thin box border, and center/top alignment.
$ws.Range("A1:N1").ClearFormats()

$ws.Cells.Item(20, 8).Value = 69989  # H20: 69994 -> 69989
$ws.Cells.Item(20, 10).Value = 69989  # J20: 69994 -> 69989
$ws.Cells.Item(20, 12).Value = 69989  # L20: 69994 -> 69989
$ws.Cells.Item(20, 14).Value = -70461  # N20: -70466 -> -70461

$ws.Cells.Item(30, 8).Value = 69989  # H30: 69994 -> 69989
$ws.Cells.Item(30, 10).Value = 69989  # J30: 69994 -> 69989
$ws.Cells.Item(30, 12).Value = 69989  # L30: 69994 -> 69989
$ws.Cells.Item(30, 14).Value = -70171  # N30: -70176 -> -70171

$ws.Cells.Item(31, 8).Value = 1265.7916  # H31: 1265.875 -> 1265.7916
$ws.Cells.Item(31, 9).Value = 1081.1428  # I31: 1081.2858 -> 1081.1428
$ws.Cells.Item(31, 11).Value = 1081.1428  # K31: 1081.2858 -> 1081.1428
$ws.Cells.Item(31, 13).Value = -786.1428000000001  # M31: -786.2858000000001 -> -786.1428000000001

$ws.Cells.Item(32, 8).Value = 3496.6667  # H32: 4625 -> 3496.6667
$ws.Cells.Item(32, 9).Value = 3496.6667  # I32: 4625 -> 3496.6667
$ws.Cells.Item(32, 11).Value = 3496.6667  # K32: 4625 -> 3496.6667
$ws.Cells.Item(32, 13).Value = -3180.6667  # M32: -4309 -> -3180.6667

$ws.Cells.Item(34, 8).Value = 1265.7916  # H34: 1265.875 -> 1265.7916
$ws.Cells.Item(34, 9).Value = 1081.1428  # I34: 1081.2858 -> 1081.1428
$ws.Cells.Item(34, 11).Value = 1081.1428  # K34: 1081.2858 -> 1081.1428
$ws.Cells.Item(34, 13).Value = -879.1428000000001  # M34: -879.2858000000001 -> -879.1428000000001

$ws.Cells.Item(128, 8).Value = 69989  # H128: 69994 -> 69989
$ws.Cells.Item(128, 10).Value = 69989  # J128: 69994 -> 69989
$ws.Cells.Item(128, 12).Value = 69989  # L128: 69994 -> 69989
$ws.Cells.Item(128, 14).Value = -79949  # N128: -79954 -> -79949

$ws.Cells.Item(129, 8).Value = 92000  # H129: 94000 -> 92000
$ws.Cells.Item(129, 10).Value = 92000  # J129: 94000 -> 92000
$ws.Cells.Item(129, 12).Value = 92000  # L129: 94000 -> 92000
$ws.Cells.Item(129, 14).Value = -102000  # N129: -104000 -> -102000

# ---------------------------------------------------------------------------
# Sheet: CUL
# ---------------------------------------------------------------------------
$ws = $wb.Worksheets.Item("CUL")

# Reset header row (A1:N1) to the default style: removes bold font,
# thin box border, and center/top alignment.
$ws.Range("A1:N1").ClearFormats()

$ws.Cells.Item(2, 8).Value = 90.47369  # H2: 102.47059 -> 90.47369
$ws.Cells.Item(2, 9).Value = 60.916668  # I2: 64.72727 -> 60.916668
$ws.Cells.Item(2, 10).Value = 141.14285  # J2: 171.66667 -> 141.14285
$ws.Cells.Item(2, 11).Value = 365.500008  # K2: 388.36362 -> 365.500008
$ws.Cells.Item(2, 12).Value = 846.8571000000001  # L2: 1030.00002 -> 846.8571000000001
$ws.Cells.Item(2, 13).Value = -252.500008  # M2: -275.36362 -> -252.500008
$ws.Cells.Item(2, 14).Value = -1072.8571  # N2: -1256.00002 -> -1072.8571

$ws.Cells.Item(34, 8).Value = 4209.8  # H34: 3892.2727 -> 4209.8
$ws.Cells.Item(34, 10).Value = 5892.7144  # J34: 5245.75 -> 5892.7144
$ws.Cells.Item(34, 12).Value = 17678.1432  # L34: 15737.25 -> 17678.1432
$ws.Cells.Item(34, 14).Value = -17846.1432  # N34: -15905.25 -> -17846.1432

$ws.Cells.Item(46, 8).Value = 1050  # H46: 1240 -> 1050
$ws.Cells.Item(46, 10).Value = 1666.6666  # J46: 2450 -> 1666.6666
$ws.Cells.Item(46, 12).Value = 4999.9998  # L46: 7350 -> 4999.9998
$ws.Cells.Item(46, 14).Value = -5181.9998  # N46: -7532 -> -5181.9998

$ws.Cells.Item(55, 8).Value = 2075.6667  # H55: 2903 -> 2075.6667
$ws.Cells.Item(55, 9).Value = 613.5  # I55: 2806.5 -> 613.5
$ws.Cells.Item(55, 10).Value = 5000  # J55: 2999.5 -> 5000
$ws.Cells.Item(55, 11).Value = 1840.5  # K55: 8419.5 -> 1840.5
$ws.Cells.Item(55, 12).Value = 15000  # L55: 8998.5 -> 15000
$ws.Cells.Item(55, 13).Value = -1663.5  # M55: -8242.5 -> -1663.5
$ws.Cells.Item(55, 14).Value = -15354  # N55: -9352.5 -> -15354

$ws.Cells.Item(68, 8).Value = 2679.7  # H68: 2693.3447 -> 2679.7
$ws.Cells.Item(68, 10).Value = 2891.0908  # J68: 2920 -> 2891.0908
$ws.Cells.Item(68, 12).Value = 8673.2724  # L68: 8760 -> 8673.2724
$ws.Cells.Item(68, 14).Value = -10295.2724  # N68: -10382 -> -10295.2724

$ws.Cells.Item(71, 8).Value = 2679.7  # H71: 2693.3447 -> 2679.7
$ws.Cells.Item(71, 10).Value = 2891.0908  # J71: 2920 -> 2891.0908
$ws.Cells.Item(71, 12).Value = 26019.8172  # L71: 26280 -> 26019.8172
$ws.Cells.Item(71, 14).Value = -34131.8172  # N71: -34392 -> -34131.8172

$ws.Cells.Item(131, 8).Value = 2413.3635  # H131: 2387.1667 -> 2413.3635
$ws.Cells.Item(131, 9).Value = 2680.7144  # I131: 2961 -> 2680.7144
$ws.Cells.Item(131, 10).Value = 1945.5  # J131: 1813.3334 -> 1945.5
$ws.Cells.Item(131, 11).Value = 8042.1432  # K131: 8883 -> 8042.1432
$ws.Cells.Item(131, 12).Value = 5836.5  # L131: 5440.0002 -> 5836.5
$ws.Cells.Item(131, 13).Value = -3002.1432  # M131: -3843 -> -3002.1432
$ws.Cells.Item(131, 14).Value = -15916.5  # N131: -15520.0002 -> -15916.5

# ---------------------------------------------------------------------------
# Sheet: GSM
# ---------------------------------------------------------------------------
$ws = $wb.Worksheets.Item("GSM")

# Reset header row (A1:N1) to the default style: removes bold font,
# thin box border, and center/top alignment.
$ws.Range("A1:N1").ClearFormats()

$ws.Cells.Item(107, 8).Value = 2269.3  # H107: 2694.875 -> 2269.3
$ws.Cells.Item(107, 9).Value = 978.1429  # I107: 1142.6 -> 978.1429
$ws.Cells.Item(107, 11).Value = 978.1429  # K107: 1142.6 -> 978.1429
$ws.Cells.Item(107, 13).Value = 941.8571  # M107: 777.4000000000001 -> 941.8571

$ws.Cells.Item(139, 8).Value = 22995  # H139: 25000 -> 22995
$ws.Cells.Item(139, 10).Value = 22995  # J139: 25000 -> 22995
$ws.Cells.Item(139, 12).Value = 22995  # L139: 25000 -> 22995
$ws.Cells.Item(139, 14).Value = -33275  # N139: -35280 -> -33275

$ws.Cells.Item(141, 8).Value = 143999.5  # H141: 153998.25 -> 143999.5
$ws.Cells.Item(141, 10).Value = 143999.5  # J141: 153998.25 -> 143999.5
$ws.Cells.Item(141, 12).Value = 143999.5  # L141: 153998.25 -> 143999.5
$ws.Cells.Item(141, 14).Value = -154359.5  # N141: -164358.25 -> -154359.5

# ---------------------------------------------------------------------------
# Sheet: LTW
# ---------------------------------------------------------------------------
$ws = $wb.Worksheets.Item("LTW")

# Reset header row (A1:N1) to the default style: removes bold font,
# thin box border, and center/top alignment.
$ws.Range("A1:N1").ClearFormats()

$ws.Cells.Item(7, 8).Value = 4392.7144  # H7: 3643.0833 -> 4392.7144
$ws.Cells.Item(7, 9).Value = 2941.25  # I7: 2748.111 -> 2941.25
$ws.Cells.Item(7, 11).Value = 2941.25  # K7: 2748.111 -> 2941.25
$ws.Cells.Item(7, 13).Value = -2829.25  # M7: -2636.111 -> -2829.25

$ws.Cells.Item(101, 8).Value = 14373.333  # H101: 19708.666 -> 14373.333
$ws.Cells.Item(101, 10).Value = 14373.333  # J101: 19708.666 -> 14373.333
$ws.Cells.Item(101, 12).Value = 14373.333  # L101: 19708.666 -> 14373.333
$ws.Cells.Item(101, 14).Value = -20863.333  # N101: -26198.666 -> -20863.333

$ws.Cells.Item(125, 8).Value = 59998  # H125: 59998.332 -> 59998
$ws.Cells.Item(125, 10).Value = 59998  # J125: 59998.332 -> 59998
$ws.Cells.Item(125, 12).Value = 59998  # L125: 59998.332 -> 59998
$ws.Cells.Item(125, 14).Value = -69838  # N125: -69838.332 -> -69838

$ws.Cells.Item(126, 8).Value = 4392.7144  # H126: 3643.0833 -> 4392.7144
$ws.Cells.Item(126, 9).Value = 2941.25  # I126: 2748.111 -> 2941.25
$ws.Cells.Item(126, 11).Value = 8823.75  # K126: 8244.332999999999 -> 8823.75
$ws.Cells.Item(126, 13).Value = -6353.75  # M126: -5774.332999999999 -> -6353.75

$ws.Cells.Item(128, 8).Value = 79997.5  # H128: 80000 -> 79997.5
$ws.Cells.Item(128, 10).Value = 79997.5  # J128: 80000 -> 79997.5
$ws.Cells.Item(128, 12).Value = 79997.5  # L128: 80000 -> 79997.5
$ws.Cells.Item(128, 14).Value = -89957.5  # N128: -89960 -> -89957.5

$ws.Cells.Item(132, 8).Value = 3955.375  # H132: 3594.111 -> 3955.375
$ws.Cells.Item(132, 9).Value = 2606.6667  # I132: 2334.8572 -> 2606.6667
$ws.Cells.Item(132, 11).Value = 7820.000100000001  # K132: 7004.571599999999 -> 7820.000100000001
$ws.Cells.Item(132, 13).Value = -5290.000100000001  # M132: -4474.571599999999 -> -5290.000100000001

$ws.Cells.Item(136, 8).Value = 38463436  # H136: 41668640 -> 38463436
$ws.Cells.Item(136, 9).Value = 1735.25  # I136: 1740.5714 -> 1735.25
$ws.Cells.Item(136, 10).Value = 100002160  # J136: 100002296 -> 100002160
$ws.Cells.Item(136, 11).Value = 5205.75  # K136: 5221.7142 -> 5205.75
$ws.Cells.Item(136, 12).Value = 300006480  # L136: 300006888 -> 300006480
$ws.Cells.Item(136, 13).Value = -2655.75  # M136: -2671.7142 -> -2655.75
$ws.Cells.Item(136, 14).Value = -300011580  # N136: -300011988 -> -300011580

# ---------------------------------------------------------------------------
# Sheet: WVR
# ---------------------------------------------------------------------------
$ws = $wb.Worksheets.Item("WVR")

# Reset header row (A1:N1) to the default style: removes bold font,
# thin box border, and center/top alignment.
$ws.Range("A1:N1").ClearFormats()

$ws.Cells.Item(26, 8).Value = 3000  # H26: 0 -> 3000
$ws.Cells.Item(26, 10).Value = 3000  # J26: 0 -> 3000
$ws.Cells.Item(26, 12).Value = 3000  # L26: 0 -> 3000
$ws.Cells.Item(26, 14).Value = -3586  # N26: None -> -3586

$ws.Cells.Item(62, 8).Value = 30376.25  # H62: 53249.5 -> 30376.25
$ws.Cells.Item(62, 9).Value = 5751  # I62: 6499 -> 5751
$ws.Cells.Item(62, 10).Value = 55001.5  # J62: 100000 -> 55001.5
$ws.Cells.Item(62, 11).Value = 5751  # K62: 6499 -> 5751
$ws.Cells.Item(62, 12).Value = 55001.5  # L62: 100000 -> 55001.5
$ws.Cells.Item(62, 13).Value = -5127  # M62: -5875 -> -5127
$ws.Cells.Item(62, 14).Value = -56249.5  # N62: -101248 -> -56249.5

$ws.Cells.Item(65, 8).Value = 30376.25  # H65: 53249.5 -> 30376.25
$ws.Cells.Item(65, 9).Value = 5751  # I65: 6499 -> 5751
$ws.Cells.Item(65, 10).Value = 55001.5  # J65: 100000 -> 55001.5
$ws.Cells.Item(65, 11).Value = 28755  # K65: 32495 -> 28755
$ws.Cells.Item(65, 12).Value = 275007.5  # L65: 500000 -> 275007.5
$ws.Cells.Item(65, 13).Value = -25635  # M65: -29375 -> -25635
$ws.Cells.Item(65, 14).Value = -281247.5  # N65: -506240 -> -281247.5

$ws.Cells.Item(69, 8).Value = 10940.7  # H69: 11863.5 -> 10940.7
$ws.Cells.Item(69, 10).Value = 10940.7  # J69: 11863.5 -> 10940.7
$ws.Cells.Item(69, 12).Value = 10940.7  # L69: 11863.5 -> 10940.7
$ws.Cells.Item(69, 14).Value = -12438.7  # N69: -13361.5 -> -12438.7

$ws.Cells.Item(72, 8).Value = 10940.7  # H72: 11863.5 -> 10940.7
$ws.Cells.Item(72, 10).Value = 10940.7  # J72: 11863.5 -> 10940.7
$ws.Cells.Item(72, 12).Value = 32822.10000000001  # L72: 35590.5 -> 32822.10000000001
$ws.Cells.Item(72, 14).Value = -40310.10000000001  # N72: -43078.5 -> -40310.10000000001

$ws.Cells.Item(98, 8).Value = 25551.5  # H98: 25554.5 -> 25551.5
$ws.Cells.Item(98, 10).Value = 25551.5  # J98: 25554.5 -> 25551.5
$ws.Cells.Item(98, 12).Value = 25551.5  # L98: 25554.5 -> 25551.5
$ws.Cells.Item(98, 14).Value = -31541.5  # N98: -31544.5 -> -31541.5

$ws.Cells.Item(107, 8).Value = 7274.875  # H107: 3144.577 -> 7274.875
$ws.Cells.Item(107, 9).Value = 7028.4287  # I107: 2769.5715 -> 7028.4287
$ws.Cells.Item(107, 10).Value = 9000  # J107: 4719.6 -> 9000
$ws.Cells.Item(107, 11).Value = 21085.2861  # K107: 8308.7145 -> 21085.2861
$ws.Cells.Item(107, 12).Value = 27000  # L107: 14158.8 -> 27000
$ws.Cells.Item(107, 13).Value = -19165.2861  # M107: -6388.7145 -> -19165.2861
$ws.Cells.Item(107, 14).Value = -30840  # N107: -17998.8 -> -30840

$ws.Cells.Item(124, 8).Value = 92000  # H124: 92500 -> 92000
$ws.Cells.Item(124, 10).Value = 92000  # J124: 92500 -> 92000
$ws.Cells.Item(124, 12).Value = 92000  # L124: 92500 -> 92000
$ws.Cells.Item(124, 14).Value = -101820  # N124: -102320 -> -101820

$ws.Cells.Item(126, 8).Value = 1993  # H126: 1995.6 -> 1993
$ws.Cells.Item(126, 9).Value = 1993.2858  # I126: 1995.6 -> 1993.2858
$ws.Cells.Item(126, 10).Value = 1991  # J126: 0 -> 1991
$ws.Cells.Item(126, 11).Value = 5979.857400000001  # K126: 5986.799999999999 -> 5979.857400000001
$ws.Cells.Item(126, 12).Value = 5973  # L126: 0 -> 5973
$ws.Cells.Item(126, 13).Value = -3509.857400000001  # M126: -3516.799999999999 -> -3509.857400000001
$ws.Cells.Item(126, 14).Value = -10913  # N126: None -> -10913

$ws.Cells.Item(132, 8).Value = 5923.3335  # H132: 6484.75 -> 5923.3335
$ws.Cells.Item(132, 9).Value = 4733.3335  # I132: 4812.1665 -> 4733.3335
$ws.Cells.Item(132, 10).Value = 8303.333  # J132: 11502.5 -> 8303.333
$ws.Cells.Item(132, 11).Value = 14200.0005  # K132: 14436.4995 -> 14200.0005
$ws.Cells.Item(132, 12).Value = 24909.999  # L132: 34507.5 -> 24909.999
$ws.Cells.Item(132, 13).Value = -11670.0005  # M132: -11906.4995 -> -11670.0005
$ws.Cells.Item(132, 14).Value = -29969.999  # N132: -39567.5 -> -29969.999

$ws.Cells.Item(140, 8).Value = 78000  # H140: 77997.5 -> 78000
$ws.Cells.Item(140, 9).Value = 0  # I140: 77995 -> 0
$ws.Cells.Item(140, 11).Value = 0  # K140: 77995 -> 0
$ws.Cells.Item(140, 13).Value = $null  # M140: -72815 -> (removed)
